$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9975135922431946
$ws.Range("B1").Value = 2.14692211151123
$ws.Range("C1").Value = 7.428821563720703
$ws.Range("D1").Value = 2.407320737838745
$ws.Range("E1").Value = 1.350888609886169
